$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 745
$ws.Range("I4").Value = 745
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 745
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -631
$ws.Range("N4").ClearContents()
# Row 11
$ws.Range("H11").Value = 7161.533
$ws.Range("I11").Value = 7161.533
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 7161.533
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -7021.533
# Row 28
$ws.Range("H28").Value = 1176.9231
$ws.Range("I28").Value = 1480.9286
$ws.Range("J28").Value = 822.25
$ws.Range("K28").Value = 1480.9286
$ws.Range("L28").Value = 822.25
$ws.Range("M28").Value = -995.9286
$ws.Range("N28").Value = -1792.25
# Row 70
$ws.Range("H70").Value = 2176
$ws.Range("I70").Value = 1400
$ws.Range("J70").Value = 2370
$ws.Range("K70").Value = 4200
$ws.Range("L70").Value = 7110
$ws.Range("M70").Value = -3930
$ws.Range("N70").Value = -7650
# Row 73
$ws.Range("H73").Value = 2176
$ws.Range("I73").Value = 1400
$ws.Range("J73").Value = 2370
$ws.Range("K73").Value = 4200
$ws.Range("L73").Value = 7110
$ws.Range("M73").Value = -3264
$ws.Range("N73").Value = -8982
# Row 92
$ws.Range("H92").Value = 1438
$ws.Range("I92").Value = 1232.8462
$ws.Range("J92").Value = 2200
$ws.Range("K92").Value = 1232.8462
$ws.Range("L92").Value = 2200
$ws.Range("M92").Value = 15.15380000000005
$ws.Range("N92").Value = -4696
# Row 111
$ws.Range("H111").Value = 6866.3335
$ws.Range("I111").Value = 6866.3335
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 20599.0005
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -17532.0005
$ws.Range("N111").ClearContents()
# Row 116
$ws.Range("H116").Value = 5573.3335
$ws.Range("I116").Value = 1431.1111
$ws.Range("J116").Value = 18000
$ws.Range("K116").Value = 1431.1111
$ws.Range("L116").Value = 18000
$ws.Range("M116").Value = 2010.8889
$ws.Range("N116").Value = -24884
# Row 127
$ws.Range("H127").Value = 1623.3529
$ws.Range("I127").Value = 748.5
$ws.Range("J127").Value = 1740
$ws.Range("K127").Value = 2245.5
$ws.Range("L127").Value = 5220
$ws.Range("M127").Value = 2714.5
$ws.Range("N127").Value = -15140
# Row 135
$ws.Range("H135").Value = 48945.617
$ws.Range("I135").Value = 56855.445
$ws.Range("J135").Value = 1486.6666
$ws.Range("K135").Value = 511699.005
$ws.Range("L135").Value = 13379.9994
$ws.Range("M135").Value = -509164.005
$ws.Range("N135").Value = -18449.9994
# Row 137
$ws.Range("H137").Value = 1819956.4
$ws.Range("I137").Value = 4763407
$ws.Range("J137").Value = 1942.6765
$ws.Range("K137").Value = 14290221
$ws.Range("L137").Value = 5828.029500000001
$ws.Range("M137").Value = -14287671
$ws.Range("N137").Value = -10928.0295

$ws = $wb.Worksheets.Item("ARM")
# Row 132
$ws.Range("H132").Value = 202279.9
$ws.Range("I132").Value = 168583.33
$ws.Range("J132").Value = 252824.75
$ws.Range("K132").Value = 505749.99
$ws.Range("L132").Value = 758474.25
$ws.Range("M132").Value = -503219.99
$ws.Range("N132").Value = -763534.25

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 634.16
$ws.Range("I31").Value = 352.42105
$ws.Range("J31").Value = 700.2469
$ws.Range("K31").Value = 352.42105
$ws.Range("L31").Value = 700.2469
$ws.Range("M31").Value = -57.42104999999998
$ws.Range("N31").Value = -1290.2469
# Row 34
$ws.Range("H34").Value = 634.16
$ws.Range("I34").Value = 352.42105
$ws.Range("J34").Value = 700.2469
$ws.Range("K34").Value = 352.42105
$ws.Range("L34").Value = 700.2469
$ws.Range("M34").Value = -150.42105
$ws.Range("N34").Value = -1104.2469
# Row 58
$ws.Range("H58").Value = 56671120
$ws.Range("I58").Value = 63754760
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 63754760
$ws.Range("L58").Value = 2000
$ws.Range("M58").Value = -63754557
$ws.Range("N58").Value = -2406
# Row 99
$ws.Range("H99").Value = 2370.2
$ws.Range("I99").Value = 2370.2
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2370.2
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -872.1999999999998
$ws.Range("N99").ClearContents()
# Row 126
$ws.Range("H126").Value = 2370.2
$ws.Range("I126").Value = 2370.2
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7110.599999999999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4640.599999999999
$ws.Range("N126").ClearContents()
# Row 132
$ws.Range("H132").Value = 205349.6
$ws.Range("I132").Value = 12750
$ws.Range("J132").Value = 253499.5
$ws.Range("K132").Value = 38250
$ws.Range("L132").Value = 760498.5
$ws.Range("M132").Value = -35720
$ws.Range("N132").Value = -765558.5
# Row 134
$ws.Range("H134").Value = 42686.54
$ws.Range("I134").Value = 937.3125
$ws.Range("J134").Value = 109485.3
$ws.Range("K134").Value = 2811.9375
$ws.Range("L134").Value = 328455.9
$ws.Range("M134").Value = -276.9375
$ws.Range("N134").Value = -333525.9
# Row 136
$ws.Range("H136").Value = 56671120
$ws.Range("I136").Value = 63754760
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 191264280
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -191261730
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("CUL")
# Row 106
$ws.Range("H106").Value = 3888.8333
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 3888.8333
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 11666.4999
$ws.Range("N106").Value = -13558.4999
# Row 107
$ws.Range("H107").Value = 1003.2958
$ws.Range("I107").Value = 683.7368
$ws.Range("J107").Value = 1120.0577
$ws.Range("K107").Value = 2051.2104
$ws.Range("L107").Value = 3360.1731
$ws.Range("M107").Value = -131.2103999999999
$ws.Range("N107").Value = -7200.1731
# Row 116
$ws.Range("H116").Value = 122049.91
$ws.Range("I116").Value = 95921.28999999999
$ws.Range("J116").Value = 167775
$ws.Range("K116").Value = 287763.87
$ws.Range("L116").Value = 503325
$ws.Range("M116").Value = -284321.87
$ws.Range("N116").Value = -510209
# Row 118
$ws.Range("H118").Value = 743
$ws.Range("I118").Value = 743
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 2229
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -986
$ws.Range("N118").ClearContents()
# Row 122
$ws.Range("H122").Value = 999.9
$ws.Range("I122").Value = 414.14285
$ws.Range("J122").Value = 2366.6667
$ws.Range("K122").Value = 3727.28565
$ws.Range("L122").Value = 21300.0003
$ws.Range("M122").Value = -1277.28565
$ws.Range("N122").Value = -26200.0003
# Row 132
$ws.Range("H132").Value = 2079.1904
$ws.Range("I132").Value = 1727.7142
$ws.Range("J132").Value = 2254.9285
$ws.Range("K132").Value = 15549.4278
$ws.Range("L132").Value = 20294.3565
$ws.Range("M132").Value = -13019.4278
$ws.Range("N132").Value = -25354.3565
# Row 134
$ws.Range("H134").Value = 2703.0908
$ws.Range("I134").Value = 1526.1111
$ws.Range("J134").Value = 7999.5
$ws.Range("K134").Value = 4578.3333
$ws.Range("L134").Value = 23998.5
$ws.Range("M134").Value = 491.6666999999998
$ws.Range("N134").Value = -34138.5
# Row 136
$ws.Range("H136").Value = 3999.3333
$ws.Range("I136").Value = 6000
$ws.Range("J136").Value = 2999
$ws.Range("K136").Value = 18000
$ws.Range("L136").Value = 8997
$ws.Range("M136").Value = -12900
$ws.Range("N136").Value = -19197
# Row 138
$ws.Range("H138").Value = 1575.9
$ws.Range("I138").Value = 1380
$ws.Range("J138").Value = 2033
$ws.Range("K138").Value = 4140
$ws.Range("L138").Value = 6099
$ws.Range("M138").Value = 1000
$ws.Range("N138").Value = -16379

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 1262.3077
$ws.Range("I93").Value = 1200.7142
$ws.Range("J93").Value = 1521
$ws.Range("K93").Value = 1200.7142
$ws.Range("L93").Value = 1521
$ws.Range("M93").Value = 47.28580000000011
$ws.Range("N93").Value = -4017
# Row 122
$ws.Range("H122").Value = 3489.4167
$ws.Range("I122").Value = 3489.4167
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10468.2501
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8018.250100000001
$ws.Range("N122").ClearContents()
# Row 132
$ws.Range("H132").Value = 37857.07
$ws.Range("I132").Value = 2797.8
$ws.Range("J132").Value = 56309.316
$ws.Range("K132").Value = 8393.400000000001
$ws.Range("L132").Value = 168927.948
$ws.Range("M132").Value = -5863.400000000001
$ws.Range("N132").Value = -173987.948

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 877
$ws.Range("I107").Value = 790
$ws.Range("J107").Value = 1225
$ws.Range("K107").Value = 2370
$ws.Range("L107").Value = 3675
$ws.Range("M107").Value = -450
$ws.Range("N107").Value = -7515
# Row 126
$ws.Range("H126").Value = 1706.5217
$ws.Range("I126").Value = 1426.4706
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 4279.4118
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -1809.4118
$ws.Range("N126").Value = -12440
# Row 132
$ws.Range("H132").Value = 54730.08
$ws.Range("I132").Value = 54445.105
$ws.Range("J132").Value = 55015.05
$ws.Range("K132").Value = 163335.315
$ws.Range("L132").Value = 165045.15
$ws.Range("M132").Value = -160805.315
$ws.Range("N132").Value = -170105.15
